# Update "想去人数" (want-to-go attendee counts) in column F across all
# sheets of the 广州-漫展信息 workbook, matching the refreshed scrape data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览"
$ws1.Range("F3").Value = 258
$ws1.Range("F4").Value = 854
$ws1.Range("F6").Value = 423
$ws1.Range("F7").Value = 624
$ws1.Range("F8").Value = 227
$ws1.Range("F10").Value = 364
$ws1.Range("F11").Value = 161
$ws1.Range("F12").Value = 723
$ws1.Range("F14").Value = 1852
$ws1.Range("F15").Value = 377
$ws1.Range("F16").Value = 4138
$ws1.Range("F17").Value = 382
$ws1.Range("F18").Value = 499
$ws1.Range("F19").Value = 12
$ws1.Range("F20").Value = 63
$ws1.Range("F21").Value = 150

# Sheet "演出"
$ws2.Range("F2").Value = 22
$ws2.Range("F7").Value = 483
$ws2.Range("F12").Value = 20
$ws2.Range("F13").Value = 99
$ws2.Range("F15").Value = 5

# Sheet "本地生活"
$ws3.Range("F2").Value = 5377
$ws3.Range("F4").Value = 296

# Sheet "全部类型"
$ws4.Range("F3").Value = 5377
$ws4.Range("F5").Value = 22
$ws4.Range("F6").Value = 296
$ws4.Range("F7").Value = 258
$ws4.Range("F12").Value = 483
$ws4.Range("F13").Value = 854
$ws4.Range("F17").Value = 423
$ws4.Range("F18").Value = 624
$ws4.Range("F19").Value = 227
$ws4.Range("F22").Value = 364
$ws4.Range("F23").Value = 161
$ws4.Range("F25").Value = 20
$ws4.Range("F26").Value = 723
$ws4.Range("F28").Value = 99
$ws4.Range("F29").Value = 1852
$ws4.Range("F30").Value = 377
$ws4.Range("F31").Value = 4138
$ws4.Range("F33").Value = 382
$ws4.Range("F34").Value = 499
$ws4.Range("F35").Value = 12
$ws4.Range("F36").Value = 63
$ws4.Range("F37").Value = 5
$ws4.Range("F38").Value = 150
